$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the alignment to A7 (top-left of the future merge) first, then merge
# A7:A8 so the merged area keeps A7's value and both cells share that style.
$a7 = $ws.Range("A7")
$a7.HorizontalAlignment = -4131  # xlLeft
$a7.VerticalAlignment = -4108    # xlCenter

$ws.Range("A7:A8").Merge()

$ws.Range("A7:A8").Select()
